$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume (E) columns to Text format so that
# numeric-looking strings (e.g. "6.260", "29.259.81") retain their exact
# textual representation instead of being auto-coerced to numbers.
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = "29.259.81"
$ws.Range("E2").Value = "  +3.01%  "
$ws.Range("D3").Value = "1.894.89"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").Value = "314.40"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "0.5144"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").Value = "0.3917"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "0.08415"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.115"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "42.24"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "6.260"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.889.09"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "20.65"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "7.300"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "93.06"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "0.00001104"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "0.06732"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "6.007"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").Value = "29.259.98"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "2.213"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "2.106.56"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "158.98"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").Value = "20.86"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").Value = "2.427"
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("D30").Value = "127.63"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.058"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.1046"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "6.152"
$ws.Range("E33").Value = "  +6.59%  "
$ws.Range("D34").Value = "3.657"
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").Value = "0.02480"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").Value = "0.06569"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "0.2191"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "9.000"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").Value = "5.176"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").Value = "1.225"
$ws.Range("E40").Value = "  +2.91%  "
$ws.Range("D41").Value = "0.6496"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "1.232"
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("D43").Value = "11.25"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "0.6048"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "13.14"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").Value = "3.673"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "2.051"
$ws.Range("E47").Value = "  +2.24%  "
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("D49").Value = "123.40"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "1.172"
$ws.Range("E50").Value = "  -1.52%  "
$ws.Range("D51").Value = "77.40"
$ws.Range("E51").Value = "  +0.66%  "

# Restore the original (default/"Normal") style on the Price/Volume cells so
# only the Value changes -- matching the source diff, which shows no style edits.
$numRng.Style = "Normal"

